# Apply the tracked changes:
#  - Planilha1!E3 value changes from 123 to 1234
#  - Planilha1 selection moves from E4 to E6, and it is no longer the
#    selected/active tab
#  - Planilha2 selection moves from C8 to B6, and it becomes the
#    selected/active tab

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Planilha2")

# Update the numeric value in Planilha1!E3.
$ws1.Range("E3").Value = 1234

# Update Planilha1's remembered selection (E4 -> E6). Activate it first so
# the selection actually "sticks" on that sheet.
$ws1.Activate()
$ws1.Range("E6").Select()

# Make Planilha2 the active tab and update its remembered selection
# (C8 -> B6). Activating it last leaves it as the workbook's active sheet
# when the file is saved.
$ws2.Activate()
$ws2.Range("B6").Select()
